$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111815519
$ws.Range("B2").Value = 77515
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = "Garnlav"
$ws.Range("G2").Value = "Alectoria sarmentosa"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("Q2").Value = 458215.7474518137
$ws.Range("R2").Value = 7054621.063481365
$ws.Range("AC2").ClearContents()

# Row 3
$ws.Range("A3").Value = 111815508
$ws.Range("B3").Value = 56398
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("Q3").Value = 458162.4570845839
$ws.Range("R3").Value = 7054329.489790585
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("AC3").Value = "ringhack"

# Row 4
$ws.Range("A4").Value = 111815512
$ws.Range("Q4").Value = 458154.6107204149
$ws.Range("R4").Value = 7054646.336103803

# Row 5
$ws.Range("A5").Value = 111815514
$ws.Range("B5").Value = 89423
$ws.Range("E5").Value = 5432
$ws.Range("F5").Value = "Granticka"
$ws.Range("G5").Value = "Porodaedalea chrysoloma"
$ws.Range("H5").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("K5").ClearContents()
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("Q5").Value = 458153.7808649908
$ws.Range("R5").Value = 7054482.19637617
$ws.Range("AC5").ClearContents()

# Row 6
$ws.Range("A6").Value = 111815517
$ws.Range("B6").Value = 77515
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("Q6").Value = 458250.8216980004
$ws.Range("R6").Value = 7054375.482693202

# Row 8
$ws.Range("A8").Value = 111815513
$ws.Range("B8").Value = 56398
$ws.Range("E8").Value = 100109
$ws.Range("F8").Value = "Tretåig hackspett"
$ws.Range("G8").Value = "Picoides tridactylus"
$ws.Range("H8").Value = "(Linnaeus, 1758)"
$ws.Range("Q8").Value = 458173.7327805056
$ws.Range("R8").Value = 7054711.474791372
$ws.Range("K8").Value = ""
$ws.Range("L8").Value = ""
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = ""
$ws.Range("AC8").Value = "ringhack gamla"

# Row 10
$ws.Range("A10").Value = 111815518
$ws.Range("Q10").Value = 458250.901553072
$ws.Range("R10").Value = 7054618.376188213

# Row 11
$ws.Range("A11").Value = 111815516
$ws.Range("B11").Value = 89423
$ws.Range("E11").Value = 5432
$ws.Range("F11").Value = "Granticka"
$ws.Range("G11").Value = "Porodaedalea chrysoloma"
$ws.Range("H11").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("K11").ClearContents()
$ws.Range("L11").ClearContents()
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()
$ws.Range("Q11").Value = 458289.5512131723
$ws.Range("R11").Value = 7054475.069158822
$ws.Range("AC11").ClearContents()

# Row 12
$ws.Range("A12").Value = 111815510
$ws.Range("Q12").Value = 458203.7272220219
$ws.Range("R12").Value = 7054385.000644128

# Row 13
$ws.Range("A13").Value = 111815507
$ws.Range("B13").Value = 56398
$ws.Range("E13").Value = 100109
$ws.Range("F13").Value = "Tretåig hackspett"
$ws.Range("G13").Value = "Picoides tridactylus"
$ws.Range("H13").Value = "(Linnaeus, 1758)"
$ws.Range("Q13").Value = 458151.5539710881
$ws.Range("R13").Value = 7054482.225765129
$ws.Range("K13").Value = ""
$ws.Range("L13").Value = ""
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = ""
$ws.Range("AC13").Value = "ringhack gamla"
